# Auto-generated Word COM-interop script
# Applies the diff: updates the date line and all 100 table-cell answers.

$d = $word.ActiveDocument

# --- Update the date paragraph above the table ---
$d.Content.Find.Execute("2025-01-17 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-18 Saturday", 2) | Out-Null

# --- Update every answer cell in the table, addressed by (row, col) ---
# (avoids relying on Find/Replace uniqueness since "73-71=2" appears twice
#  in the source document with two different replacement values)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "77-22=55"
$t.Cell(1, 2).Range.Text = "2+33=35"
$t.Cell(1, 3).Range.Text = "83-32=51"
$t.Cell(1, 4).Range.Text = "80-9=71"
$t.Cell(1, 5).Range.Text = "90-15=75"
$t.Cell(2, 1).Range.Text = "5+29=34"
$t.Cell(2, 2).Range.Text = "33-30=3"
$t.Cell(2, 3).Range.Text = "18+53=71"
$t.Cell(2, 4).Range.Text = "97-8=89"
$t.Cell(2, 5).Range.Text = "0+42=42"
$t.Cell(3, 1).Range.Text = "92-21=71"
$t.Cell(3, 2).Range.Text = "88-81=7"
$t.Cell(3, 3).Range.Text = "72-16=56"
$t.Cell(3, 4).Range.Text = "1+34=35"
$t.Cell(3, 5).Range.Text = "83-36=47"
$t.Cell(4, 1).Range.Text = "3+39=42"
$t.Cell(4, 2).Range.Text = "22+2=24"
$t.Cell(4, 3).Range.Text = "71-32=39"
$t.Cell(4, 4).Range.Text = "20+48=68"
$t.Cell(4, 5).Range.Text = "78-65=13"
$t.Cell(5, 1).Range.Text = "6+38=44"
$t.Cell(5, 2).Range.Text = "66-49=17"
$t.Cell(5, 3).Range.Text = "72+23=95"
$t.Cell(5, 4).Range.Text = "32+41=73"
$t.Cell(5, 5).Range.Text = "81-35=46"
$t.Cell(6, 1).Range.Text = "11+72=83"
$t.Cell(6, 2).Range.Text = "23-12=11"
$t.Cell(6, 3).Range.Text = "25+51=76"
$t.Cell(6, 4).Range.Text = "87-5=82"
$t.Cell(6, 5).Range.Text = "76-26=50"
$t.Cell(7, 1).Range.Text = "47-31=16"
$t.Cell(7, 2).Range.Text = "74-57=17"
$t.Cell(7, 3).Range.Text = "27-10=17"
$t.Cell(7, 4).Range.Text = "62-40=22"
$t.Cell(7, 5).Range.Text = "45+33=78"
$t.Cell(8, 1).Range.Text = "78-7=71"
$t.Cell(8, 2).Range.Text = "25+39=64"
$t.Cell(8, 3).Range.Text = "2-1=1"
$t.Cell(8, 4).Range.Text = "9+82=91"
$t.Cell(8, 5).Range.Text = "89-76=13"
$t.Cell(9, 1).Range.Text = "59+38=97"
$t.Cell(9, 2).Range.Text = "7+1=8"
$t.Cell(9, 3).Range.Text = "50-15=35"
$t.Cell(9, 4).Range.Text = "8+42=50"
$t.Cell(9, 5).Range.Text = "52-26=26"
$t.Cell(10, 1).Range.Text = "24+3=27"
$t.Cell(10, 2).Range.Text = "52-5=47"
$t.Cell(10, 3).Range.Text = "52-7=45"
$t.Cell(10, 4).Range.Text = "93-0=93"
$t.Cell(10, 5).Range.Text = "0+15=15"
$t.Cell(11, 1).Range.Text = "50-24=26"
$t.Cell(11, 2).Range.Text = "58-0=58"
$t.Cell(11, 3).Range.Text = "48+34=82"
$t.Cell(11, 4).Range.Text = "90-29=61"
$t.Cell(11, 5).Range.Text = "2+7=9"
$t.Cell(12, 1).Range.Text = "57+14=71"
$t.Cell(12, 2).Range.Text = "74-23=51"
$t.Cell(12, 3).Range.Text = "63+3=66"
$t.Cell(12, 4).Range.Text = "39+0=39"
$t.Cell(12, 5).Range.Text = "22+73=95"
$t.Cell(13, 1).Range.Text = "83-78=5"
$t.Cell(13, 2).Range.Text = "1+91=92"
$t.Cell(13, 3).Range.Text = "32+11=43"
$t.Cell(13, 4).Range.Text = "8+11=19"
$t.Cell(13, 5).Range.Text = "3+13=16"
$t.Cell(14, 1).Range.Text = "78+16=94"
$t.Cell(14, 2).Range.Text = "26-4=22"
$t.Cell(14, 3).Range.Text = "70+24=94"
$t.Cell(14, 4).Range.Text = "99-5=94"
$t.Cell(14, 5).Range.Text = "27+64=91"
$t.Cell(15, 1).Range.Text = "49+0=49"
$t.Cell(15, 2).Range.Text = "64-30=34"
$t.Cell(15, 3).Range.Text = "21-13=8"
$t.Cell(15, 4).Range.Text = "45+12=57"
$t.Cell(15, 5).Range.Text = "19+25=44"
$t.Cell(16, 1).Range.Text = "72+22=94"
$t.Cell(16, 2).Range.Text = "1+65=66"
$t.Cell(16, 3).Range.Text = "96-58=38"
$t.Cell(16, 4).Range.Text = "35+57=92"
$t.Cell(16, 5).Range.Text = "79-1=78"
$t.Cell(17, 1).Range.Text = "6+76=82"
$t.Cell(17, 2).Range.Text = "91-13=78"
$t.Cell(17, 3).Range.Text = "23-3=20"
$t.Cell(17, 4).Range.Text = "82-30=52"
$t.Cell(17, 5).Range.Text = "79-52=27"
$t.Cell(18, 1).Range.Text = "53-27=26"
$t.Cell(18, 2).Range.Text = "43-21=22"
$t.Cell(18, 3).Range.Text = "64-19=45"
$t.Cell(18, 4).Range.Text = "16+41=57"
$t.Cell(18, 5).Range.Text = "77-71=6"
$t.Cell(19, 1).Range.Text = "59-37=22"
$t.Cell(19, 2).Range.Text = "35+45=80"
$t.Cell(19, 3).Range.Text = "66-15=51"
$t.Cell(19, 4).Range.Text = "2+30=32"
$t.Cell(19, 5).Range.Text = "55+27=82"
$t.Cell(20, 1).Range.Text = "31+38=69"
$t.Cell(20, 2).Range.Text = "77-20=57"
$t.Cell(20, 3).Range.Text = "52+3=55"
$t.Cell(20, 4).Range.Text = "23+41=64"
$t.Cell(20, 5).Range.Text = "43+26=69"
